# Reorder the language rows by descending value and drop the last two
# (Uzbek, Vietnamese) rows so only 20 languages remain (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English", 26.67548697766628),
    @("Spanish", 7.972132198076175),
    @("Japanese", 7.861116838679116),
    @("Chinese", 7.006145155861569),
    @("German", 6.341494018219962),
    @("Arabic", 4.666235999545704),
    @("Portuguese", 3.91425188231456),
    @("French", 3.676450157697122),
    @("Italian", 3.57573665938741),
    @("Russian", 3.484647200064367),
    @("Malay-Indonesian", 2.839055281566525),
    @("Dutch", 1.648062347215147),
    @("Persian", 1.428440225274951),
    @("Korean", 1.422469181731608),
    @("Turkish", 1.359196299608752),
    @("Thai", 1.054986118876026),
    @("Polish", 0.8251759685676816),
    @("Urdu", 0.7962262383900328),
    @("Swedish", 0.5212865891345118),
    @("Bengali", 0.4075418315839745)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-obsolete rows 22 and 23 (previously Uzbek / Vietnamese)
$ws.Range("A22:B23").Delete()
